$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D").Insert()

# Copy number formats from column E (which now holds the old column D formatting) into new column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# Populate new column D with the latest period values
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 349600
$ws.Range("D9").Value = 160100
$ws.Range("D10").Value = 189500
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 340400
$ws.Range("D18").Value = 9200
$ws.Range("D20").Value = -1300
$ws.Range("D21").Value = 14700
$ws.Range("D22").Value = 1000
$ws.Range("D23").Value = 6900
$ws.Range("D24").Value = 600
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 6200
$ws.Range("D27").Value = 6200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 1300
$ws.Range("D33").Value = 6200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 6200
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 11900
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 41300
$ws.Range("D44").Value = 52900
$ws.Range("D45").Value = 800
$ws.Range("D46").Value = 106900
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 23500
$ws.Range("D49").Value = 33200
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 33500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 197100
$ws.Range("D57").Value = 15200
$ws.Range("D58").Value = 12000
$ws.Range("D59").Value = 39000
$ws.Range("D60").Value = 66200
$ws.Range("D61").Value = 5200
$ws.Range("D62").Value = 26500
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 98000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 77300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 99200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 6200
$ws.Range("D83").Value = 6900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 20300
$ws.Range("D91").Value = -2500
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -7800
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -4500
$ws.Range("D101").Value = -500
$ws.Range("D102").Value = 7500
